# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E2) and
# "Correspond Handback DateTime" (H2) timestamps on the zh-cn and de-de
# sheets to reflect the latest handback run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-21 17:09:57"
$wsZhCn.Range("H2").Value = "2016-03-21 17:10:27"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-21 17:10:01"
$wsDeDe.Range("H2").Value = "2016-03-21 17:10:34"
